$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.300.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.684.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.41%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.00"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5229"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.72%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2703"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06402"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.98"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07510"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.704.78"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.560"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5786"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008426"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.20"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.349.39"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.913"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.86"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.58"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.178"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.684"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1229"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.80"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06628"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +12.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.343"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.71%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.567"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.564"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.655"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.025"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6188"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.401"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.695"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.376"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.106.28"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01613"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8755"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.78"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.832.72"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.65%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.63"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.131"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05270"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4309"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.037"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.92%  "
